# Apply the "Fragen" questionnaire-string refactor to Tabelle1 (sheet1).
# The shared strings describing each question's field grammar were
# reworded (the ":pflicht" suffix style became an "|pflicht" infix style,
# a couple of entries picked up proper field labels/prefixes, etc.) and
# the sheet's cells were repointed at the new text accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D2").Value  = "Name/Firma:text|pflicht;Anschrift:text|pflicht;Steuernummer:text|pflicht;AMA-Betriebsnummer:text;SVS-Versicherungsnummer:text;Bankverbindung(IBAN/BIC):text;Familienstand:select(ledig,verh.,geschieden,verwitwet)|pflicht;Kinder:number;Vollmacht-und-DSGVO-Einwilligung:checkbox:pflicht"
$ws.Range("D3").Value  = "Erwerbsart:select(Vollerwerb, Nebenerwerb)|pflicht"
$ws.Range("D4").Value  = "HINWEIS --> AMA Daten werden von der Finanz verwertet, saubere Erstdaten helfen bei Plausibilitätsprüfung:info;Flächenaufstellung:select(Eigenbewirtschaftung,Pacht,Mitbewirtschaftung)|pflicht;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D5").Value  = "Tierarten-und-Bestände:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D6").Value  = "Anzahl Hektar:number;typische Nutzung:text|pflicht;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("C7").Value  = "A3. Nebentätigkeiten gesamt (Erstaufnahme)"
$ws.Range("D15").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D16").Value = "Art:select(Neubau,Umbau,Zubau);Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D17").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D19").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D20").Value = "Förderungen:mcheckbox(Keine,AMA,ÖPUL,Sonstige);Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D21").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D22").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D23").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D24").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"

$ws.Range("D25").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
# D25 previously had no explicit cell style; the edited workbook gives it
# the same "Text" number format (style index 1) used by its neighbours.
$ws.Range("D25").NumberFormat = "@"

$ws.Range("D26").Value = "Richtig-und-Vollständig:checkbox|pflicht"
$ws.Range("C27").Value = "Sie haben das Fragebogenende erreicht. Beenden Sie mit JA."
$ws.Range("D28").Value = "Name/Firma:text|pflicht;Anschrift:text|pflicht;Steuernummer:text;AMA-Betriebsnummer:text;SVS-Versicherungsnummer:text;Bankverbindung(IBAN/BIC)|pflicht:text;Familienstand:select(ledig,verh.,geschieden,verwitwet);Kinder:number;Eingabe-vollständig:checkbox|pflicht"
$ws.Range("D29").Value = "Betriebsführer:text;Mitunternehmer:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D30").Value = "Erwerbsart:select(Vollerwerb,Nebenerwerb);Bewirtschaftungsart:select(Bio,konventionell);Kurzbeschreibung:text;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D31").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D32").Value = "Beschreibung:text|pflicht;Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"
$ws.Range("D33").Value = "Richtig-und-Vollständig:checkbox|pflicht"
$ws.Range("C34").Value = "Sie haben das Fragebogenende erreicht. Beenden Sie mit JA."

# Update the sheet view: scroll so row 7 is at the top and select D34
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("D34").Select()
